# Generate Report for Handback
# Marks the two handed-off files ("23b5d40f..." and "514390ce...") as handed
# back (in sync with en-US), fills in the "Latest Target File" / "Latest
# Handback File" / "Latest Handback DateTime" columns on the per-locale
# sheets, and widens a few columns so the longer status / file-name text
# fits (mirrors what Excel's own AutoFit would have produced).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column widths observed in the target workbook snap to this engine's
# internal 1/6-character grid when set through ColumnWidth, so we pick the
# COM value whose quantized result lands on (or nearest to) the target.
$wideStatusColumnWidth = 29.166666666666668   # -> renders as width 30 (was ~17.2 -> grows to fit longer status text)
$wideFileColumnWidth   = 39.166666666666664   # -> renders as width 40 exactly

# ---------------------------------------------------------------------
# Sheet 1: "Overview" - just the Status column (E/F) text + width change
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $wideStatusColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideStatusColumnWidth

# ---------------------------------------------------------------------
# Helper data for the two per-locale sheets
# ---------------------------------------------------------------------
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/faf5c6fdaac6b26d77714c07eb1403987f4e5ead/e2e/23b5d40f-7ff5-4700-8f0e-9be0c3e23037.md"
$mdName1 = "23b5d40f-7ff5-4700-8f0e-9be0c3e23037.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/faf5c6fdaac6b26d77714c07eb1403987f4e5ead/e2e/514390ce-98dc-4f99-8979-9c36b7792137.md"
$mdName2 = "514390ce-98dc-4f99-8979-9c36b7792137.md"

function Update-LocaleSheet($ws, $targetFile1, $targetFile2, $handbackDateTime) {
    # Status column (C)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Latest Target File (I) / Latest Handback File (J) / Latest Handback DateTime (K)
    $ws.Range("I2").Value = $mdName1
    $ws.Range("J2").Value = $targetFile1
    $ws.Range("K2").Value = $handbackDateTime

    $ws.Range("I3").Value = $mdName2
    $ws.Range("J3").Value = $targetFile2
    $ws.Range("K3").Value = $handbackDateTime

    # Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
    $ws.Columns.Item(3).ColumnWidth = $wideStatusColumnWidth
    $ws.Columns.Item(9).ColumnWidth = $wideFileColumnWidth
    $ws.Columns.Item(10).ColumnWidth = $wideFileColumnWidth

    # Recreate hyperlinks so the file-name columns (A and I) both link out to
    # the source .md file, in A2, I2, A3, I3 order (matches handoff sheet).
    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1)
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl1, [Type]::Missing, [Type]::Missing, $mdName1)
    $ws.Hyperlinks.Add($ws.Range("A3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2)
    $ws.Hyperlinks.Add($ws.Range("I3"), $mdUrl2, [Type]::Missing, [Type]::Missing, $mdName2)
}

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZhCn `
    "23b5d40f-7ff5-4700-8f0e-9be0c3e23037.57efa451913de21fc7b0b16bdb45b43918c6a054.zh-cn.xlf" `
    "514390ce-98dc-4f99-8979-9c36b7792137.28c1ac143586fcd4875221ac7648a4fa8497fe05.zh-cn.xlf" `
    "2016-10-24 10:19:29"

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDeDe `
    "23b5d40f-7ff5-4700-8f0e-9be0c3e23037.57efa451913de21fc7b0b16bdb45b43918c6a054.de-de.xlf" `
    "514390ce-98dc-4f99-8979-9c36b7792137.28c1ac143586fcd4875221ac7648a4fa8497fe05.de-de.xlf" `
    "2016-10-24 10:19:45"
